# Update the dSF column (F) values to reflect the repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    3  = -3
    6  = -9
    7  = -8
    8  = 11
    9  = -11
    13 = -6
    16 = -7
    18 = -3
    20 = 0
    22 = 4
    23 = -5
    25 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
